$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A71").Value = "me20233253@sva.edu.eg"
$ws.Range("B71").Value = "https://mariam2005856.github.io/mariam_essam/"

$ws.Range("A72").Value = "re20231348@sva.edu.eg"
$ws.Range("B72").Value = "https://rehabezzat122.github.io/first-web-page/"

$ws.Range("A73").Value = "ea20230350@sva.edu.eg"
$ws.Range("B73").Value = "https://eman2005820.github.io/mypage./"
